$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.188.18"
$ws.Range("E2").Value = "  +2.03%  "
$ws.Range("D3").Value = "2.780.17"
$ws.Range("E3").Value = "  +2.99%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.28"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.55"
$ws.Range("E6").Value = "  +8.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.621"
$ws.Range("E7").Value = "  +2.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.998"
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("D9").Value = "2.793.99"
$ws.Range("E9").Value = "  +2.51%  "
$ws.Range("E10").Value = "  +1.51%  "
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.401"
$ws.Range("E12").Value = "  +3.46%  "
$ws.Range("E13").Value = "  +1.30%  "
$ws.Range("D14").Value = "3.274.33"
$ws.Range("E14").Value = "  +3.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.68"
$ws.Range("E15").Value = "  +4.68%  "
$ws.Range("D16").Value = "64.087.57"
$ws.Range("E16").Value = "  +2.14%  "
$ws.Range("E17").Value = "  +6.04%  "
$ws.Range("D18").Value = "2.787.87"
$ws.Range("E18").Value = "  +2.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.50"
$ws.Range("E19").Value = "  +4.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.05"
$ws.Range("E20").Value = "  +3.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "368.19"
$ws.Range("E22").Value = "  +1.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.577"
$ws.Range("E23").Value = "  +8.53%  "
$ws.Range("E24").Value = "  +0.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.56"
$ws.Range("E25").Value = "  +3.69%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.176"
$ws.Range("E26").Value = "  +6.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.79"
$ws.Range("E27").Value = "  +2.38%  "
$ws.Range("D28").Value = "0.0₃0970"
$ws.Range("E28").Value = "  +13.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("E30").Value = "  +1.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.29"
$ws.Range("E31").Value = "  +2.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.28"
$ws.Range("E32").Value = "  +9.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "172.57"
$ws.Range("E33").Value = "  +1.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.12"
$ws.Range("E34").Value = "  +7.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "20.85"
$ws.Range("E35").Value = "  +1.71%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  +5.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.83"
$ws.Range("E38").Value = "  +0.64%  "
$ws.Range("E39").Value = "  +1.63%  "
$ws.Range("E40").Value = "  +1.50%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.32"
$ws.Range("E41").Value = "  +11.21%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "342.85"
$ws.Range("E42").Value = "  -1.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.97"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.51"
$ws.Range("E44").Value = "  +4.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.80"
$ws.Range("E45").Value = "  +5.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0614"
$ws.Range("E46").Value = "  +3.56%  "
$ws.Range("E47").Value = "  +2.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0263"
$ws.Range("E48").Value = "  +1.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "138.92"
$ws.Range("E49").Value = "  +1.40%  "
$ws.Range("E50").Value = "  +2.26%  "
$ws.Range("D51").Value = "2.181.99"
$ws.Range("E51").Value = "  +2.26%  "
